# edit.ps1 -- apply the "Updated cryptos list" data refresh to cryptos.xlsx
#
# The sheet stores every data cell (Coin / Link / Price / Volume(1h)) as plain
# text, even when the text looks like a number (e.g. "213.31", "0.492").
# If we just did $range.Value = "213.31", Excel would happily reinterpret it
# as the number 213.31 and the round trip would lose the original text-cell
# typing/formatting. So for any new value that parses as a number we first
# force the cell to Text format, assign the value, and then restore the cell
# style to "Normal" so we do not leave a stray number-format style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Price cells whose new text would otherwise be auto-detected as a number ---
Set-TextValue "D5" '213.31'
Set-TextValue "D7" '0.492'
Set-TextValue "D10" '19.08'
Set-TextValue "D14" '4.05'
Set-TextValue "D15" '0.523'
Set-TextValue "D17" '62.98'
Set-TextValue "D19" '209.99'
Set-TextValue "D21" '4.31'
Set-TextValue "D22" '9.40'
Set-TextValue "D23" '6.13'
Set-TextValue "D24" '1.92'
Set-TextValue "D25" '146.24'
Set-TextValue "D28" '15.43'
Set-TextValue "D29" '6.70'
Set-TextValue "D30" '0.0515'
Set-TextValue "D32" '3.22'
Set-TextValue "D45" '92.29'
Set-TextValue "D46" '1.55'
Set-TextValue "D48" '54.60'
Set-TextValue "D49" '0.0513'
Set-TextValue "D50" '0.410'
Set-TextValue "D51" '7.54'

# --- Remaining cells (coin names, links, "Volume(1h)" percentages, and Price
#     values that are naturally non-numeric text such as "26.653.35") ---
$ws.Range("D2").Value = '26.653.35'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '1.636.57'
$ws.Range("E3").Value = '  +1.55%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("E7").Value = '  +1.07%  '
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("E9").Value = '  +0.88%  '
$ws.Range("E10").Value = '  +3.42%  '
$ws.Range("E11").Value = '  +2.63%  '
$ws.Range("D12").Value = '1.866.09'
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D13").Value = '1.634.54'
$ws.Range("E13").Value = '  +1.20%  '
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").Value = '26.663.53'
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").Value = '0.0₃0741'
$ws.Range("E18").Value = '  +1.54%  '
$ws.Range("E19").Value = '  +3.36%  '
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("E24").Value = '  +2.74%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  -1.38%  '
$ws.Range("E28").Value = '  +1.35%  '
$ws.Range("E29").Value = '  +1.78%  '
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("E33").Value = '  +0.63%  '
$ws.Range("E34").Value = '  +0.76%  '
$ws.Range("E35").Value = '  -0.49%  '
$ws.Range("D36").Value = '1.167.44'
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("E38").Value = '  +2.21%  '
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("E42").Value = '  +1.53%  '
$ws.Range("E43").Value = '  +1.91%  '
$ws.Range("D44").Value = '1.774.64'
$ws.Range("E44").Value = '  +1.34%  '
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("E46").Value = '  +1.53%  '
$ws.Range("E47").Value = '  +9.83%  '
$ws.Range("E48").Value = '  +0.45%  '
$ws.Range("E49").Value = '  +1.01%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("E50").Value = '  +0.56%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E51").Value = '  +3.95%  '
